# teste 2 anos.xlsx - apply the authored edits
#
# Summary of the change (from the OOXML diff):
#   - Row 2 "Item" answer-key cells are updated:
#       AK2 (Item2) : "X" -> cleared (blank)
#       AN2 (Item5) : "-" -> "x" (new lowercase marker, distinct shared string)
#       AO2 (Item6) : "X" -> "x"
#       AP2 (Item7) : "-" -> cleared (blank)
#   - The active selection on the sheet moves from AT9 to AS2.
#
# (The workbook.xml absPath/documentId GUID churn and the styles.xml dxf
# renumbering visible in the raw XML diff are incidental artifacts that
# Excel/engines regenerate on every save - they carry no business meaning
# and are not driven explicitly here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the two cells that lost their marker.
$ws.Range("AK2").Value = ""
$ws.Range("AP2").Value = ""

# Replace the marker text on the other two cells with the lowercase "x".
$ws.Range("AN2").Value = "x"
$ws.Range("AO2").Value = "x"

# Move the active selection to AS2 (was AT9).
$ws.Range("AS2").Select() | Out-Null
